$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Tnfsf12", "Cd163", "ECs", 3, 1, 5.526037, 16.578111, 0.1907740523064932, 0.1907740523064932, 2, 0.6666666666666666, 0.1082676666666667, 0.324803, 0.02641792511238088, 0.02641792511238088, 0.5982911319036667, 5.384620187133, 0.005039854627218371, 0.005039854627218372),
    @("ECs", "Tnfsf12", "Cd163", "M2", 3, 1, 5.526037, 16.578111, 0.1907740523064932, 0.1907740523064932, 3, 1, 3.049833, 9.149499, 0.7441765605545629, 0.7441765605545632, 16.853490001821, 151.681410016389, 0.1419695780885024, 0.1419695780885025),
    @("ECs", "Tnfsf12", "Cd163", "sCs", 3, 1, 5.526037, 16.578111, 0.1907740523064932, 0.1907740523064932, 3, 1, 0.9401646666666666, 2.820494, 0.229405514333056, 0.229405514333056, 5.195384734092666, 46.758462606834, 0.04376461959077241, 0.04376461959077241),
    @("FAPs", "Tnfsf12", "Cd163", "ECs", 3, 1, 8.502533, 25.507599, 0.2935309110814287, 0.2935309110814287, 2, 0.6666666666666666, 0.1082676666666667, 0.324803, 0.02641792511238088, 0.02641792511238088, 0.9205494086663333, 8.284944677997, 0.007754477627118114, 0.007754477627118114),
    @("FAPs", "Tnfsf12", "Cd163", "M2", 3, 1, 8.502533, 25.507599, 0.2935309110814287, 0.2935309110814287, 3, 1, 3.049833, 9.149499, 0.7441765605545629, 0.7441765605545632, 25.931305726989, 233.381751542901, 0.2184388238250249, 0.2184388238250249),
    @("FAPs", "Tnfsf12", "Cd163", "sCs", 3, 1, 8.502533, 25.507599, 0.2935309110814287, 0.2935309110814287, 3, 1, 0.9401646666666666, 2.820494, 0.229405514333056, 0.229405514333056, 7.993781103767333, 71.94402993390599, 0.06733760962928569, 0.06733760962928569),
    @("M2", "Tnfsf12", "Cd163", "ECs", 3, 1, 10.61443666666667, 31.84331, 0.3664396557335078, 0.3664396557335078, 2, 0.6666666666666666, 0.1082676666666667, 0.324803, 0.02641792511238088, 0.02641792511238088, 1.149200290881111, 10.34280261793, 0.00968057538337444, 0.009680575383374442),
    @("M2", "Tnfsf12", "Cd163", "M2", 3, 1, 10.61443666666667, 31.84331, 0.3664396557335078, 0.3664396557335078, 3, 1, 3.049833, 9.149499, 0.7441765605545629, 0.7441765605545632, 32.37225922241, 291.35033300169, 0.27269580265456, 0.2726958026545601),
    @("M2", "Tnfsf12", "Cd163", "sCs", 3, 1, 10.61443666666667, 31.84331, 0.3664396557335078, 0.3664396557335078, 3, 1, 0.9401646666666666, 2.820494, 0.229405514333056, 0.229405514333056, 9.979318310571111, 89.81386479513999, 0.08406327769557334, 0.08406327769557335),
    @("sCs", "Tnfsf12", "Cd163", "ECs", 3, 1, 4.323390666666667, 12.970172, 0.1492553808785702, 0.1492553808785702, 2, 0.6666666666666666, 0.1082676666666667, 0.324803, 0.02641792511238088, 0.02641792511238088, 0.4680834195684445, 4.212750776116001, 0.003943017474669954, 0.003943017474669953),
    @("sCs", "Tnfsf12", "Cd163", "M2", 3, 1, 4.323390666666667, 12.970172, 0.1492553808785702, 0.1492553808785702, 3, 1, 3.049833, 9.149499, 0.7441765605545629, 0.7441765605545632, 13.185619527092, 118.670575743828, 0.1110723559864757, 0.1110723559864757),
    @("sCs", "Tnfsf12", "Cd163", "sCs", 3, 1, 4.323390666666667, 12.970172, 0.1492553808785702, 0.1492553808785702, 3, 1, 0.9401646666666666, 2.820494, 0.229405514333056, 0.229405514333056, 4.064699144996445, 36.58229230496801, 0.03424000741742458, 0.03424000741742458)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $rowVals = $rows[$i]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value2 = $rowVals[$c]
    }
}

Write-Host "done"
